$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''27.666.84'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.05%  '
$ws.Range("D3").Value = '''1.873.18'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.83%  '
$ws.Range("E4").Value = '  +0.24%  '
$ws.Range("D5").Value = '''331.81'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.60%  '
$ws.Range("E6").Value = '  +0.21%  '
$ws.Range("D7").Value = '''0.4718'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +4.06%  '
$ws.Range("D8").Value = '''0.3942'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.85%  '
$ws.Range("D9").Value = '''47.89'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.89%  '
$ws.Range("D10").Value = '''0.08042'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.41%  '
$ws.Range("D11").Value = '''1.027'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.97%  '
$ws.Range("D12").Value = '''22.06'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.08%  '
$ws.Range("D13").Value = '''1.891.72'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.44%  '
$ws.Range("D14").Value = '''5.959'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.66%  '
$ws.Range("D15").Value = '''7.120'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.12%  '
$ws.Range("D16").Value = '''1.004'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.25%  '
$ws.Range("D17").Value = '''0.00001048'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.53%  '
$ws.Range("D18").Value = '''86.88'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.03%  '
$ws.Range("D19").Value = '''0.06679'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.05%  '
$ws.Range("D20").Value = '''17.14'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.64%  '
$ws.Range("D21").Value = '''1.004'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.21%  '
$ws.Range("D22").Value = '''27.679.50'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.10%  '
$ws.Range("D23").Value = '''5.514'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.32%  '
$ws.Range("D24").Value = '''10.96'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.06%  '
$ws.Range("D25").Value = '''2.309'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.21%  '
$ws.Range("D26").Value = '''2.096.56'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.54%  '
$ws.Range("D27").Value = '''158.66'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.10%  '
$ws.Range("D28").Value = '''20.13'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.54%  '
$ws.Range("E29").Value = '  +1.14%  '
$ws.Range("D30").Value = '''5.562'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.17%  '
$ws.Range("D31").Value = '''122.24'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.03%  '
$ws.Range("D32").Value = '''0.9732'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.76%  '
$ws.Range("E33").Value = '  +2.50%  '
$ws.Range("D34").Value = '''1.446'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.69%  '
$ws.Range("D35").Value = '''3.594'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.07%  '
$ws.Range("D36").Value = '''5.330'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.33%  '
$ws.Range("D37").Value = '''0.06104'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.87%  '
$ws.Range("E38").Value = '  +0.58%  '
$ws.Range("E39").Value = '  +0.27%  '
$ws.Range("D40").Value = '''8.137'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.75%  '
$ws.Range("D41").Value = '''0.6019'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.77%  '
$ws.Range("D42").Value = '''0.1901'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.35%  '
$ws.Range("E43").Value = '  +0.98%  '
$ws.Range("D44").Value = '''1.266'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.74%  '
$ws.Range("D45").Value = '''0.5693'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.41%  '
$ws.Range("D46").Value = '''12.13'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.68%  '
$ws.Range("E47").Value = '  +0.86%  '
$ws.Range("D48").Value = '''3.380'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.15%  '
$ws.Range("D49").Value = '''0.06887'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.78%  '
$ws.Range("D50").Value = '''114.79'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.87%  '
$ws.Range("E51").Value = '  +10.27%  '
